$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.03140789642930031
$ws.Range("C2").Value = 0.016744256019592285
$ws.Range("D2").Value = 0.00962611474096775
$ws.Range("E2").Value = 0.004523232579231262
$ws.Range("F2").Value = 0.0005770409479737282
$ws.Range("G2").Value = 0.0007732643280178308
$ws.Range("J2").Value = 0.12755978107452393
$ws.Range("K2").Value = 1.4510115385055542
